$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("table_placard_bas")

# Update data values in row 2 (hauteur / etagere)
$ws.Range("B2").Value = 1
$ws.Range("E2").Value = 2

# Update the active selection, as recorded in the saved view state
$ws.Activate()
$ws.Range("C7").Select()
